# Staging.IndicatorLocation.xlsx — rename the staging key columns.
#
# Sheet1 row 2 holds three "for internal use" placeholder labels in A2:C2:
#   A2: IndicatorLocation_ID  -> IndicatorBusinessKey
#   B2: IndicatorSourceKey    -> IndicatorLocation_ID
#   C2: LocationSourceKey     -> LocationBusinessKey
#
# (This mirrors the repo rename where *SourceKey columns were renamed to
# *BusinessKey, and IndicatorLocation_ID moved from the A column into B.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "IndicatorBusinessKey"
$ws.Range("B2").Value = "IndicatorLocation_ID"
$ws.Range("C2").Value = "LocationBusinessKey"
